# Apply the cryptos list update (price/volume refresh + two coin swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store numeric-looking values as plain
# text (e.g. "69.134.21", "0.492", "1.00"). Force a Text number format across the
# whole data range first so that COM does not silently coerce these strings into
# numbers (which would drop formatting like trailing zeros or thousands-style dots).
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '69.134.21'
$ws.Range("D3").Value = '3.678.30'
$ws.Range("E3").Value = '  -2.62%  '
$ws.Range("D5").Value = '679.62'
$ws.Range("E5").Value = '  -3.31%  '
$ws.Range("D6").Value = '161.65'
$ws.Range("E6").Value = '  -4.01%  '
$ws.Range("D7").Value = '3.674.85'
$ws.Range("E7").Value = '  -2.82%  '
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  -5.00%  '
$ws.Range("E10").Value = '  -7.20%  '
$ws.Range("D11").Value = '7.20'
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("D12").Value = '0.448'
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("D13").Value = '0.0000235'
$ws.Range("E13").Value = '  -6.55%  '
$ws.Range("D14").Value = '33.32'
$ws.Range("E14").Value = '  -7.14%  '
$ws.Range("D15").Value = '4.303.10'
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").Value = '3.671.03'
$ws.Range("E16").Value = '  -1.37%  '
$ws.Range("D17").Value = '69.219.28'
$ws.Range("E17").Value = '  -1.59%  '
$ws.Range("E18").Value = '  -1.71%  '
$ws.Range("D19").Value = '16.24'
$ws.Range("E19").Value = '  -5.45%  '
$ws.Range("D20").Value = '6.59'
$ws.Range("E20").Value = '  -7.10%  '
$ws.Range("D21").Value = '483.30'
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("D22").Value = '9.77'
$ws.Range("E22").Value = '  -7.05%  '
$ws.Range("D23").Value = '0.660'
$ws.Range("E23").Value = '  -8.39%  '
$ws.Range("D24").Value = '79.40'
$ws.Range("E24").Value = '  -6.24%  '
$ws.Range("D25").Value = '3.832.35'
$ws.Range("E25").Value = '  -2.39%  '
$ws.Range("E26").Value = '  -10.95%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '11.49'
$ws.Range("E28").Value = '  -4.11%  '
$ws.Range("D29").Value = '9.46'
$ws.Range("E29").Value = '  -8.80%  '
$ws.Range("E30").Value = '  -10.60%  '
$ws.Range("D31").Value = '2.73'
$ws.Range("E31").Value = '  -10.82%  '
$ws.Range("D32").Value = '2.10'
$ws.Range("E32").Value = '  -4.87%  '
$ws.Range("D33").Value = '6.69'
$ws.Range("E33").Value = '  -7.91%  '
$ws.Range("D34").Value = '1.00'
$ws.Range("D35").Value = '26.73'
$ws.Range("E35").Value = '  -7.65%  '
$ws.Range("E36").Value = '  -5.26%  '
$ws.Range("D37").Value = '3.648.56'
$ws.Range("E37").Value = '  -2.54%  '
$ws.Range("D38").Value = '8.47'
$ws.Range("E38").Value = '  -5.56%  '
$ws.Range("D39").Value = '6.06'
$ws.Range("E39").Value = '  +3.15%  '
$ws.Range("D40").Value = '0.0932'
$ws.Range("E40").Value = '  -7.35%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = '2.17'
$ws.Range("E42").Value = '  -4.96%  '
$ws.Range("D44").Value = '0.954'
$ws.Range("E44").Value = '  -7.77%  '
$ws.Range("D45").Value = '158.51'
$ws.Range("E45").Value = '  -3.31%  '
$ws.Range("D46").Value = '47.95'
$ws.Range("E46").Value = '  -1.60%  '
$ws.Range("D47").Value = '2.80'
$ws.Range("E47").Value = '  -13.62%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value = '0.000276'
$ws.Range("E48").Value = '  -9.75%  '
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").Value = '389.28'
$ws.Range("E49").Value = '  -6.49%  '
$ws.Range("E50").Value = '  -5.34%  '
$ws.Range("D51").Value = '8.01'
$ws.Range("E51").Value = '  -7.22%  '

# Restore the plain default style on the data range (no explicit style / number
# format should remain attached to these cells, matching the source workbook).
$priceVolRange.Style = "Normal"

Write-Output "Applied cryptos list update"
